$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - first tax record (was blank placeholder row)
$ws.Range("A2").Value = "004/ZZZ"
$ws.Range("B2").Value = "Direction régionale"
$ws.Range("C2").Value = "IR801997"
$ws.Range("D2").Value = "NOUBAIL MOHAMMED"
$ws.Range("E2").Value = "non"
$ws.Range("F2").Value = "mensuelle"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 1000
$ws.Range("O2").Value = 1000

# Row 3 - second tax record (new row)
$ws.Range("A3").Value = "004/ZZZ"
$ws.Range("B3").Value = "Direction régionale"
$ws.Range("C3").Value = "IB19558"
$ws.Range("D3").Value = "ZERNAKH ABDELLAH"
$ws.Range("E3").Value = "non"
$ws.Range("F3").Value = "mensuelle"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 1000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 1000

# Row 4 - totals row (new, blank text cells + summed numbers)
$ws.Range("A4").Value = " "
$ws.Range("B4").Value = " "
$ws.Range("C4").Value = " "
$ws.Range("D4").Value = " "
$ws.Range("E4").Value = " "
$ws.Range("F4").Value = " "
$ws.Range("G4").Value = " "
$ws.Range("H4").Value = 2000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 2000
